$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 2693
$ws1.Range("F7").Value = 2218
$ws1.Range("F8").Value = 1813
$ws1.Range("F9").Value = 214
$ws1.Range("F11").Value = 2469
$ws1.Range("F12").Value = 543
$ws1.Range("F13").Value = 232
$ws1.Range("F16").Value = 124
$ws1.Range("F17").Value = 111
$ws1.Range("F18").Value = 9161
$ws1.Range("F20").Value = 7104
$ws1.Range("F21").Value = 11621
$ws1.Range("F24").Value = 233
$ws1.Range("F25").Value = 344
$ws1.Range("F26").Value = 555
$ws1.Range("F27").Value = 2562
$ws1.Range("I27").Value = "//i0.hdslb.com/bfs/openplatform/202410/4K3QpDLU1728885983390.jpeg"
$ws1.Range("F30").Value = 2502
$ws1.Range("F31").Value = 686
$ws1.Range("F33").Value = 4504
$ws1.Range("F34").Value = 872
$ws1.Range("F35").Value = 351
$ws1.Range("F37").Value = 521
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F11").Value = 6
$ws2.Range("F14").Value = 66
$ws2.Range("F16").Value = 99
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 146
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 2693
$ws4.Range("F9").Value = 2218
$ws4.Range("F11").Value = 1813
$ws4.Range("F13").Value = 214
$ws4.Range("F14").Value = 2469
$ws4.Range("F16").Value = 543
$ws4.Range("F17").Value = 232
$ws4.Range("F20").Value = 124
$ws4.Range("F21").Value = 111
$ws4.Range("F22").Value = 9161
$ws4.Range("F24").Value = 7104
$ws4.Range("F25").Value = 11621
$ws4.Range("F28").Value = 233
$ws4.Range("F29").Value = 344
$ws4.Range("F31").Value = 555
$ws4.Range("F32").Value = 6
$ws4.Range("F33").Value = 2562
$ws4.Range("I33").Value = "//i0.hdslb.com/bfs/openplatform/202410/4K3QpDLU1728885983390.jpeg"
$ws4.Range("F39").Value = 4504
$ws4.Range("F40").Value = 66
$ws4.Range("F42").Value = 99
$ws4.Range("F46").Value = 521
